$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.11820000000002
$ws.Range("C4").Value = -10.85969999999999
$ws.Range("D4").Value = -6.453399999999999

$ws.Range("C5").Value = -14.54380000000002

$ws.Range("A7").Value = -21.6133

$ws.Range("C8").Value = -11.8207

$ws.Range("D9").Value = -7.749600000000003

$ws.Range("A16").Value = -20.16889999999999
$ws.Range("C16").Value = -12.15479999999999

$ws.Range("D18").Value = -8.351799999999992
